$ErrorActionPreference = "Stop"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: remove "Updated on" timestamp (F5) and timezone (I5) ---
$ws.Range("F5").Value = ""
$ws.Range("I5").Value = ""

# --- Row 23: clear the stray "ARB" label (style is kept) ---
$ws.Range("B23").Value = ""

# --- Row 24 (CCC6 block): date label change ---
$ws.Range("J24").Value = "18-May"

# --- Row 32 (APCC block): clear date label ---
$ws.Range("J32").Value = ""

# --- Rows 34-37 (APCC Desktop/HYBRID1/HYBRID2/Server): clear start/end time text ---
foreach ($r in 34..37) {
    $ws.Range("C$r").Value = ""
    $ws.Range("D$r").Value = ""
    $ws.Range("F$r").Value = ""
    $ws.Range("G$r").Value = ""
}

# --- Row 40 (ICC block): clear date label ---
$ws.Range("J40").Value = ""

# --- Rows 42-44 (ICC Frontend lines): clear start/end time text ---
foreach ($r in 42..44) {
    $ws.Range("C$r").Value = ""
    $ws.Range("D$r").Value = ""
    $ws.Range("F$r").Value = ""
    $ws.Range("G$r").Value = ""
}

# --- Rows 45-47 (ICC Backend lines): clear start/end time text ---
foreach ($r in 45..47) {
    $ws.Range("C$r").Value = ""
    $ws.Range("D$r").Value = ""
}

# --- Row 48 (EMFP block): date label change ---
$ws.Range("J48").Value = "18-May"

# --- Row 56 (BRH1 block): clear date label ---
$ws.Range("J56").Value = ""

# --- Rows 58-61 (BRH1 Notebook/Desktop/Server/AIO): clear formulas/values
#     and drop the custom [hh]:mm:ss number format back to General ---
$ws.Range("E9").Copy() | Out-Null
$ws.Range("C58:D61").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F58:G61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
foreach ($r in 58..61) {
    $ws.Range("C$r").Value = ""
    $ws.Range("D$r").Value = ""
    $ws.Range("E$r").Value = ""
    $ws.Range("F$r").Value = ""
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Value = ""
}

# --- New rows 65-66: Overtime e-mail table for EMFP ---
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B65").PasteSpecial(-4122) | Out-Null
$ws.Range("C65").PasteSpecial(-4122) | Out-Null
$ws.Range("M8").Copy() | Out-Null
$ws.Range("B66").PasteSpecial(-4122) | Out-Null
$ws.Range("C66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B65").Value = "Factory Overtime"
$ws.Range("C65").Value = "E-Mail content"
$ws.Range("B66").Value = "EMFP"

$emailText = @'
Dear All,
Please be informed that on Monday, the 28th of February, EMFP morning shift B will work overtime in the following pattern:
CSG LOB                  6:05 am -2.05pm
ISG LOB                   6:00 am -2.00pm
SHIPPING               6:05 am  -2.05pm     
Please provide relevant support.
Please note that overtime may be canceled for unpredicted, important reasons.
Manager on duty during overtime;
Regards
Mariusz Kaczewiak
EMFP Supervisor, Production Operations
Mobile: +48 500 216 562
Dell EMFP Poland
mariusz.kaczewiak@dell.com <mailto:mariusz.kaczewiak@dell.com> 
Dell Products (Poland) Sp. z o.o, Łódź, ul. Informatyczna 1
&
Maciej Kamiński
EMFP Manufacturing Operations Supervisor 
Dell | EMFP Materials
mobile +48 500 216 519
maciej_kaminski@dell.com <mailto:maciej_kaminski@dell.com> 
Dell Products (Poland) Sp. z o.o.,ul. Informatyczna 1

'@
$ws.Range("C66").Value = $emailText

Write-Host "done"